$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K3").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("K4").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("K5").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("K6").Value = "6,09 TL - 12,19 TL - 152,35 TL"
$ws.Range("K8").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("K9").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("K10").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("K11").Value = "3,05 TL - 6,09 TL - 76,17 TL"
$ws.Range("K12").Value = "WU: ,USD–; Diğer: 404,16 TL–3.403,42 TL"
$ws.Range("C13").Value = "Hesaba: Asgari 0 TL | Azami 9.999.999.999.999 TL"
$ws.Range("E13").Value = "Hesaba: Asgari 1 TL | Azami 851,5 TL"
$ws.Range("K13").Value = "Hesaba: Asgari 1 TL | Azami 865,75 TL"
$ws.Range("K14").Value = "914,14 TL - 4.265,98 TL"
